$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "ok" marker cells that appear in the updated diagram
$ws.Range("L2").Value = "ok"
$ws.Range("H13").Value = "ok"
$ws.Range("G18").Value = "ok"
$ws.Range("K18").Value = "ok"

# Restaurant menu page class: rename two operations, drop two, add/replace others
$ws.Range("K10").Value = "       +     View restaurant data ( )"
$ws.Range("K11").Value = "       +     View item data ( )"
$ws.Range("K12").Value = "       +     Place order ( )"
$ws.Range("K13").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("K15").Value = ""

# Confirmation page class: replace order name/quantity rows, add confirm order,
# move delivery address / loyalty points entries
$ws.Range("K28").Value = "       +     View delivery address ( )"
$ws.Range("K29").Value = "       +     confirm order ( )"
$ws.Range("I30").Value = "       +     View Loyalty points ( )"
$ws.Range("K30").Value = "       +     use Loyalty points ( )"
$ws.Range("K31").Value = ""
$ws.Range("K32").Value = ""
$ws.Range("K33").Value = ""
$ws.Range("K34").Value = ""

# Selection cursor moved to G11 (matches the saved view state in the diff)
$ws.Range("G11").Select()
